# Update countries & provincias Spain
# Applies the COVID-19 dashboard data refresh described in the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 24 de Abril de 2020 a las 09:52"

# Rusia (row 13)
$ws.Range("B13").Value = 68622
$ws.Range("C13").Value = 5849
$ws.Range("D13").Value = 5568
$ws.Range("E13").Value = 62439
$ws.Range("G13").Value = 60
$ws.Range("H13").Value = 615

# Uzbekistan (row 68)
$ws.Range("D68").Value = 582
$ws.Range("E68").Value = 1189

# Armenia moves ahead of Estonia/Azerbaiyan in the sorted country list.
# Row 71 now shows Armenia's refreshed figures, while the data previously
# held by rows 71/72 (Estonia, Azerbaiyan) shifts down one row.
$ws.Range("A71").Value = "Armenia"
$ws.Range("B71").Value = 1596
$ws.Range("C71").Value = 73
$ws.Range("D71").Value = 728
$ws.Range("E71").Value = 841
$ws.Range("F71").Value = 10
$ws.Range("G71").Value = 3
$ws.Range("H71").Value = 27

$ws.Range("A72").Value = "Estonia"
$ws.Range("B72").Value = 1592
$ws.Range("C72").Value = 0
$ws.Range("D72").Value = 192
$ws.Range("E72").Value = 1355
$ws.Range("F72").Value = 7
$ws.Range("G72").Value = 0
$ws.Range("H72").Value = 45

$ws.Range("A73").Value = "Azerbaiyan"
$ws.Range("B73").Value = 1548
$ws.Range("C73").Value = 0
$ws.Range("D73").Value = 948
$ws.Range("E73").Value = 580
$ws.Range("F73").Value = 14
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 20

# Letonia (row 92)
$ws.Range("B92").Value = 784
$ws.Range("C92").Value = 6
$ws.Range("E92").Value = 639
$ws.Range("G92").Value = 1
$ws.Range("H92").Value = 12

# Montenegro (row 119)
$ws.Range("B119").Value = 319
$ws.Range("C119").Value = 3
$ws.Range("E119").Value = 191

# El Salvador (row 125)
$ws.Range("D125").Value = 72
$ws.Range("E125").Value = 181

# Fiyi (row 184)
$ws.Range("D184").Value = 10
$ws.Range("E184").Value = 8
